$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = "Caroline"
    3  = "Niels"
    4  = "Ivan "
    5  = "Viktor"
    6  = "Geraldine"
    7  = "Nasrin"
    8  = "Karel"
    9  = "Alfiya"
    10 = "Nathalie"
    11 = "Yanina"
    12 = "Daryoush"
    13 = "Ariana"
    14 = "Em"
    15 = "Fabienne"
    16 = "Mahsa"
    17 = "Danil"
    18 = "Sweta"
    19 = "Gerrit"
    20 = "Alexander"
    21 = "Alice"
    22 = "Afaf"
    23 = "Jens"
    24 = "Miguel"
    25 = "Andrea"
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}
